$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values as text (e.g. "1.00", "66.140.53") which must
# stay text so Excel does not coerce them into numbers and strip formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.140.53'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '3.489.32'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '601.83'
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").Value = '143.11'
$ws.Range("E6").Value = '  -2.93%  '

$ws.Range("D7").Value = '3.492.77'
$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").Value = '8.14'
$ws.Range("E10").Value = '  +6.96%  '

$ws.Range("D11").Value = '0.135'
$ws.Range("E11").Value = '  -4.26%  '

$ws.Range("D12").Value = '0.413'
$ws.Range("E12").Value = '  -2.25%  '

$ws.Range("D13").Value = '4.086.07'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").Value = '0.0000203'
$ws.Range("E14").Value = '  -4.28%  '

$ws.Range("D15").Value = '30.23'
$ws.Range("E15").Value = '  -3.43%  '

$ws.Range("D16").Value = '3.503.11'
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '66.223.02'
$ws.Range("E17").Value = '  -1.01%  '

$ws.Range("E18").Value = '  -0.38%  '

$ws.Range("D19").Value = '10.48'
$ws.Range("E19").Value = '  +4.61%  '

$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  -3.65%  '

$ws.Range("D21").Value = '14.77'
$ws.Range("E21").Value = '  -3.27%  '

$ws.Range("D22").Value = '419.94'
$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("E23").Value = '  -2.71%  '

$ws.Range("D24").Value = '77.70'
$ws.Range("E24").Value = '  -1.50%  '

$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").Value = '0.0000116'
$ws.Range("E26").Value = '  -3.55%  '

$ws.Range("E27").Value = '  -3.24%  '

$ws.Range("D28").Value = '7.98'
$ws.Range("E28").Value = '  -4.46%  '

$ws.Range("D29").Value = '2.46'
$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.38%  '

$ws.Range("D31").Value = '0.162'
$ws.Range("E31").Value = '  -2.53%  '

$ws.Range("D32").Value = '1.48'
$ws.Range("E32").Value = '  -6.31%  '

$ws.Range("D33").Value = '25.19'
$ws.Range("E33").Value = '  -0.23%  '

$ws.Range("D34").Value = '3.488.45'
$ws.Range("E34").Value = '  +0.71%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").Value = '1.71'
$ws.Range("E36").Value = '  -4.25%  '

$ws.Range("D37").Value = '5.57'
$ws.Range("E37").Value = '  -6.00%  '

$ws.Range("D38").Value = '7.65'
$ws.Range("E38").Value = '  -2.81%  '

$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").Value = '169.52'
$ws.Range("E40").Value = '  -2.56%  '

$ws.Range("D41").Value = '0.0868'
$ws.Range("E41").Value = '  -1.30%  '

$ws.Range("D42").Value = '0.893'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").Value = '5.12'
$ws.Range("E43").Value = '  -5.00%  '

$ws.Range("D44").Value = '1.92'
$ws.Range("E44").Value = '  -8.77%  '

$ws.Range("D45").Value = '45.57'
$ws.Range("E45").Value = '  -1.51%  '

$ws.Range("D46").Value = '26.23'
$ws.Range("E46").Value = '  -9.25%  '

$ws.Range("D47").Value = '1.21'
$ws.Range("E47").Value = '  -2.12%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '7.11'
$ws.Range("E48").Value = '  -4.21%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.34'
$ws.Range("E49").Value = '  -3.37%  '

$ws.Range("D50").Value = '0.934'
$ws.Range("E50").Value = '  -4.13%  '

$ws.Range("D51").Value = '0.236'
$ws.Range("E51").Value = '  -3.13%  '
